$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = -8.824099999999993
$ws.Range("A12").Value = -22.79870000000002
$ws.Range("D12").Value = -8.292400000000002
$ws.Range("D14").Value = -8.704400000000001
$ws.Range("D22").Value = -7.9779
